$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Copy formatting (cell styles) from column C to column D for all 115 rows in one shot.
#    This also creates blank D cells with correct styles (matching B/C empty separator rows).
$ws.Range("C1:C115").Copy() | Out-Null
$ws.Range("D1:D115").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 2) D1 holds a date-formatted label ("2025/10/30") that must stay literal text, not be
#    auto-converted to a date serial by value-parsing heuristics. Enter it with a leading
#    apostrophe to force text, then restore the original column-C cell style on top.
$ws.Range("D1").Value = "'2025/10/30"
$ws.Range("C1").Copy() | Out-Null
$ws.Range("D1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# 3) Populate the remaining column D cells with their source values.
$ws.Cells.Item(2, 4).Value = "上证"
$ws.Cells.Item(3, 4).Value = 64.75
$ws.Cells.Item(4, 4).Value = 4018.86
$ws.Cells.Item(6, 4).Value = 50.68
$ws.Cells.Item(7, 4).Value = 5706.14
$ws.Cells.Item(9, 4).Value = 56.62
$ws.Cells.Item(10, 4).Value = 4754.15
$ws.Cells.Item(12, 4).Value = 59.23
$ws.Cells.Item(13, 4).Value = 7457.57
$ws.Cells.Item(15, 4).Value = 30.15
$ws.Cells.Item(16, 4).Value = 2764.78
$ws.Cells.Item(18, 4).Value = 97.16
$ws.Cells.Item(19, 4).Value = 6890.59
$ws.Cells.Item(21, 4).Value = 69.89
$ws.Cells.Item(22, 4).Value = 84750.89999999999
$ws.Cells.Item(24, 4).Value = 85.89
$ws.Cells.Item(25, 4).Value = 19909.14
$ws.Cells.Item(27, 4).Value = 79.73999999999999
$ws.Cells.Item(28, 4).Value = 39894.54
$ws.Cells.Item(30, 4).Value = 57.7
$ws.Cells.Item(31, 4).Value = 5662.42
$ws.Cells.Item(33, 4).Value = 9.779999999999999
$ws.Cells.Item(34, 4).Value = 33753.61
$ws.Cells.Item(36, 4).Value = 32.94
$ws.Cells.Item(37, 4).Value = 3533.61
$ws.Cells.Item(39, 4).Value = 51.54
$ws.Cells.Item(40, 4).Value = 3316.64
$ws.Cells.Item(42, 4).Value = 19.3
$ws.Cells.Item(43, 4).Value = 7406.51
$ws.Cells.Item(45, 4).Value = 29.98
$ws.Cells.Item(46, 4).Value = 8946.77
$ws.Cells.Item(48, 4).Value = 10.27
$ws.Cells.Item(49, 4).Value = 13108.84
$ws.Cells.Item(51, 4).Value = 26.07
$ws.Cells.Item(52, 4).Value = 12740.21
$ws.Cells.Item(54, 4).Value = 17.84
$ws.Cells.Item(55, 4).Value = 9661.34
$ws.Cells.Item(57, 4).Value = 24.54
$ws.Cells.Item(58, 4).Value = 16052.78
$ws.Cells.Item(60, 4).Value = 32.08
$ws.Cells.Item(61, 4).Value = 17526.85
$ws.Cells.Item(63, 4).Value = 21.85
$ws.Cells.Item(64, 4).Value = 10691.52
$ws.Cells.Item(66, 4).Value = 15.52
$ws.Cells.Item(67, 4).Value = 9778.9
$ws.Cells.Item(69, 4).Value = 21.86
$ws.Cells.Item(70, 4).Value = 3282.42
$ws.Cells.Item(72, 4).Value = 45.56
$ws.Cells.Item(73, 4).Value = 6109.78
$ws.Cells.Item(75, 4).Value = 27.83
$ws.Cells.Item(76, 4).Value = 9408.66
$ws.Cells.Item(78, 4).Value = 18.93
$ws.Cells.Item(79, 4).Value = 2504.5
$ws.Cells.Item(81, 4).Value = 55.61
$ws.Cells.Item(82, 4).Value = 2961.7
$ws.Cells.Item(84, 4).Value = 58.79
$ws.Cells.Item(85, 4).Value = 3030.02
$ws.Cells.Item(87, 4).Value = 52.77
$ws.Cells.Item(88, 4).Value = 4127.88
$ws.Cells.Item(90, 4).Value = 45.41
$ws.Cells.Item(91, 4).Value = 2047.92
$ws.Cells.Item(93, 4).Value = 27.82
$ws.Cells.Item(94, 4).Value = 13755.79
$ws.Cells.Item(96, 4).Value = 87.84999999999999
$ws.Cells.Item(97, 4).Value = 9687.620000000001
$ws.Cells.Item(99, 4).Value = 56.79
$ws.Cells.Item(100, 4).Value = 12602.02
$ws.Cells.Item(102, 4).Value = 4.4
$ws.Cells.Item(103, 4).Value = 2258.48
$ws.Cells.Item(105, 4).Value = 30.64
$ws.Cells.Item(106, 4).Value = 903.78
$ws.Cells.Item(108, 4).Value = 29.68
$ws.Cells.Item(109, 4).Value = 2880.9
$ws.Cells.Item(111, 4).Value = 20.55
$ws.Cells.Item(112, 4).Value = 4055.81
$ws.Cells.Item(114, 4).Value = 29.02
$ws.Cells.Item(115, 4).Value = 3402.18
